$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "43.440.89"
$ws.Range("E2").Value = "  -1.94%  "
$ws.Range("D3").Value = "2.255.11"
$ws.Range("E3").Value = "  -0.54%  "
$ws.Range("E4").Value = "  -0.16%  "
$ws.Range("D5").Value = "'231.66"
$ws.Range("E5").Value = "  +0.47%  "
$ws.Range("D6").Value = "'0.640"
$ws.Range("E6").Value = "  +0.97%  "
$ws.Range("D7").Value = "'64.32"
$ws.Range("E7").Value = "  +0.96%  "
$ws.Range("E8").Value = "  -0.17%  "
$ws.Range("D9").Value = "'0.438"
$ws.Range("E9").Value = "  -1.97%  "
$ws.Range("D10").Value = "'0.0954"
$ws.Range("E10").Value = "  -8.27%  "
$ws.Range("D11").Value = "'56.86"
$ws.Range("E11").Value = "  -0.29%  "
$ws.Range("D12").Value = "'26.39"
$ws.Range("E12").Value = "  +0.18%  "
$ws.Range("E13").Value = "  -1.15%  "
$ws.Range("D14").Value = "2.590.15"
$ws.Range("E14").Value = "  -0.52%  "
$ws.Range("D15").Value = "'14.96"
$ws.Range("E15").Value = "  -4.95%  "
$ws.Range("D16").Value = "'6.05"
$ws.Range("E16").Value = "  -1.40%  "
$ws.Range("D17").Value = "'0.823"
$ws.Range("E17").Value = "  -1.71%  "
$ws.Range("D18").Value = "2.249.83"
$ws.Range("E18").Value = "  -0.88%  "
$ws.Range("D19").Value = "43.347.96"
$ws.Range("E19").Value = "  -1.73%  "
$ws.Range("D20").Value = "0.0₃0968"
$ws.Range("E20").Value = "  -4.88%  "
$ws.Range("D21").Value = "'72.92"
$ws.Range("E21").Value = "  -0.79%  "
$ws.Range("D22").Value = "'6.08"
$ws.Range("E22").Value = "  +0.80%  "
$ws.Range("D23").Value = "'247.11"
$ws.Range("E23").Value = "  -1.81%  "
$ws.Range("D24").Value = "'3.92"
$ws.Range("E24").Value = "  +17.37%  "
$ws.Range("D25").Value = "'1.00"
$ws.Range("E25").Value = "  -0.03%  "
$ws.Range("D26").Value = "'2.44"
$ws.Range("E26").Value = "  +0.35%  "
$ws.Range("E27").Value = "  -2.25%  "
$ws.Range("B28").Value = "Cosmos"
$ws.Range("C28").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D28").Value = "'9.70"
$ws.Range("E28").Value = "  -3.01%  "
$ws.Range("B29").Value = "Monero"
$ws.Range("C29").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D29").Value = "'173.81"
$ws.Range("E29").Value = "  +0.96%  "
$ws.Range("D30").Value = "'21.63"
$ws.Range("E30").Value = "  +4.02%  "
$ws.Range("E31").Value = "  +3.88%  "
$ws.Range("D32").Value = "'0.130"
$ws.Range("E32").Value = "  -5.06%  "
$ws.Range("D33").Value = "'0.125"
$ws.Range("E33").Value = "  +0.51%  "
$ws.Range("D34").Value = "'4.93"
$ws.Range("E34").Value = "  +4.07%  "
$ws.Range("D35").Value = "'0.0680"
$ws.Range("E35").Value = "  -0.45%  "
$ws.Range("D36").Value = "'4.90"
$ws.Range("E36").Value = "  +0.76%  "
$ws.Range("D37").Value = "'3.63"
$ws.Range("E37").Value = "  -4.84%  "
$ws.Range("D38").Value = "'6.41"
$ws.Range("E38").Value = "  -3.82%  "
$ws.Range("D39").Value = "'2.27"
$ws.Range("E39").Value = "  -1.65%  "
$ws.Range("D40").Value = "'0.0250"
$ws.Range("E40").Value = "  -3.42%  "
$ws.Range("D41").Value = "'0.998"
$ws.Range("E41").Value = "  -0.38%  "
$ws.Range("D42").Value = "'8.81"
$ws.Range("E42").Value = "  +6.18%  "
$ws.Range("D43").Value = "'4.50"
$ws.Range("E43").Value = "  +3.36%  "
$ws.Range("D44").Value = "'17.14"
$ws.Range("E44").Value = "  -1.75%  "
$ws.Range("D45").Value = "'96.82"
$ws.Range("E45").Value = "  -0.83%  "
$ws.Range("D46").Value = "'0.0941"
$ws.Range("E46").Value = "  -2.42%  "
$ws.Range("E47").Value = "  -0.83%  "
$ws.Range("D48").Value = "'10.11"
$ws.Range("E48").Value = "  +2.98%  "
$ws.Range("B49").Value = "TerraClassic"
$ws.Range("C49").Value = "https://coinranking.com/coin/AaQUAs2Mc+terraclassic-lunc"
$ws.Range("D49").Value = "'0.000206"
$ws.Range("E49").Value = "  -1.90%  "
$ws.Range("B50").Value = "Maker"
$ws.Range("C50").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D50").Value = "1.428.58"
$ws.Range("E50").Value = "  -0.87%  "
$ws.Range("D51").Value = "'2.26"
$ws.Range("E51").Value = "  -1.66%  "
